$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5068.5
$ws.Range("I51").Value = 4478
$ws.Range("K51").Value = 4478
$ws.Range("M51").Value = -3994
$ws.Range("H94").Value = 2277
$ws.Range("I94").Value = 1888
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 1888
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -1437
$ws.Range("N94").Value = -5902
$ws.Range("H98").Value = 3031.2222
$ws.Range("J98").Value = 9000
$ws.Range("L98").Value = 9000
$ws.Range("N98").Value = -11996
$ws.Range("H106").Value = 2303.1538
$ws.Range("I106").Value = 2303.1538
$ws.Range("K106").Value = 2303.1538
$ws.Range("M106").Value = -1672.1538
$ws.Range("H122").Value = 3031.2222
$ws.Range("J122").Value = 9000
$ws.Range("L122").Value = 27000
$ws.Range("N122").Value = -31900
$ws.Range("H125").Value = 1015.6667
$ws.Range("J125").Value = 920.1429
$ws.Range("L125").Value = 8281.286100000001
$ws.Range("N125").Value = -13201.2861

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3308.5386
$ws.Range("I32").Value = 2145.4424
$ws.Range("J32").Value = 7960.923
$ws.Range("K32").Value = 2145.4424
$ws.Range("L32").Value = 7960.923
$ws.Range("M32").Value = -1858.4424
$ws.Range("N32").Value = -8534.922999999999
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8415.533
$ws.Range("I134").Value = 10536.947
$ws.Range("K134").Value = 31610.841
$ws.Range("M134").Value = -29075.841

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2975.7083
$ws.Range("I31").Value = 1354.0588
$ws.Range("J31").Value = 6914
$ws.Range("K31").Value = 1354.0588
$ws.Range("L31").Value = 6914
$ws.Range("M31").Value = -1059.0588
$ws.Range("N31").Value = -7504
$ws.Range("H34").Value = 2975.7083
$ws.Range("I34").Value = 1354.0588
$ws.Range("J34").Value = 6914
$ws.Range("K34").Value = 1354.0588
$ws.Range("L34").Value = 6914
$ws.Range("M34").Value = -1152.0588
$ws.Range("N34").Value = -7318
$ws.Range("H58").Value = 1351.8572
$ws.Range("J58").Value = 1629.3636
$ws.Range("L58").Value = 1629.3636
$ws.Range("N58").Value = -2035.3636
$ws.Range("H94").Value = 1082.2
$ws.Range("J94").Value = 1082.2
$ws.Range("L94").Value = 1082.2
$ws.Range("N94").Value = -1984.2
$ws.Range("H95").Value = 25204.666
$ws.Range("J95").Value = 25204.666
$ws.Range("L95").Value = 25204.666
$ws.Range("N95").Value = -30696.666
$ws.Range("H122").Value = 1477.1111
$ws.Range("I122").Value = 1562.8235
$ws.Range("J122").Value = 1331.4
$ws.Range("K122").Value = 4688.470499999999
$ws.Range("L122").Value = 3994.2
$ws.Range("M122").Value = -2238.470499999999
$ws.Range("N122").Value = -8894.2
$ws.Range("H136").Value = 1351.8572
$ws.Range("J136").Value = 1629.3636
$ws.Range("L136").Value = 4888.0908
$ws.Range("N136").Value = -9988.0908

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2766.6667
$ws.Range("J39").Value = 2766.6667
$ws.Range("L39").Value = 8300.000100000001
$ws.Range("N39").Value = -8888.000100000001
$ws.Range("H68").Value = 999.5
$ws.Range("J68").Value = 999.5
$ws.Range("L68").Value = 2998.5
$ws.Range("N68").Value = -4620.5
$ws.Range("H71").Value = 999.5
$ws.Range("J71").Value = 999.5
$ws.Range("L71").Value = 8995.5
$ws.Range("N71").Value = -17107.5
$ws.Range("H131").Value = 5690527.5
$ws.Range("I131").Value = 83333920
$ws.Range("J131").Value = 9303.951
$ws.Range("K131").Value = 250001760
$ws.Range("L131").Value = 27911.853
$ws.Range("M131").Value = -249996720
$ws.Range("N131").Value = -37991.853

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2385.158
$ws.Range("I80").Value = 2189.3333
$ws.Range("J80").Value = 2475.5386
$ws.Range("K80").Value = 2189.3333
$ws.Range("L80").Value = 2475.5386
$ws.Range("M80").Value = -1191.3333
$ws.Range("N80").Value = -4471.5386
$ws.Range("H83").Value = 2385.158
$ws.Range("I83").Value = 2189.3333
$ws.Range("J83").Value = 2475.5386
$ws.Range("K83").Value = 10946.6665
$ws.Range("L83").Value = 12377.693
$ws.Range("M83").Value = -5954.666499999999
$ws.Range("N83").Value = -22361.693
$ws.Range("H126").Value = 79777.92
$ws.Range("J126").Value = 202120.2
$ws.Range("L126").Value = 606360.6000000001
$ws.Range("N126").Value = -611300.6000000001
$ws.Range("H141").Value = 22991.125
$ws.Range("J141").Value = 22991.125
$ws.Range("L141").Value = 22991.125
$ws.Range("N141").Value = -33351.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6324.5835
$ws.Range("I7").Value = 3400
$ws.Range("J7").Value = 8413.571
$ws.Range("K7").Value = 3400
$ws.Range("L7").Value = 8413.571
$ws.Range("M7").Value = -3288
$ws.Range("N7").Value = -8637.571
$ws.Range("H40").Value = 7419.8423
$ws.Range("I40").Value = 3381.182
$ws.Range("K40").Value = 3381.182
$ws.Range("M40").Value = -3245.182
$ws.Range("H46").Value = 1737.375
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H126").Value = 6324.5835
$ws.Range("I126").Value = 3400
$ws.Range("J126").Value = 8413.571
$ws.Range("K126").Value = 10200
$ws.Range("L126").Value = 25240.713
$ws.Range("M126").Value = -7730
$ws.Range("N126").Value = -30180.713
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1708335.4
$ws.Range("I3").Value = 5002500
$ws.Range("K3").Value = 5002500
$ws.Range("M3").Value = -5002386
$ws.Range("H108").Value = 53599.6
$ws.Range("J108").Value = 53599.6
$ws.Range("L108").Value = 53599.6
$ws.Range("N108").Value = -61279.6
$ws.Range("H122").Value = 57422
$ws.Range("I122").Value = 112645.86
$ws.Range("J122").Value = 2198.1428
$ws.Range("K122").Value = 337937.58
$ws.Range("L122").Value = 6594.428400000001
$ws.Range("M122").Value = -335487.58
$ws.Range("N122").Value = -11494.4284
$ws.Range("H126").Value = 5307.36
$ws.Range("I126").Value = 4478.1577
$ws.Range("J126").Value = 7933.1665
$ws.Range("K126").Value = 13434.4731
$ws.Range("L126").Value = 23799.4995
$ws.Range("M126").Value = -10964.4731
$ws.Range("N126").Value = -28739.4995
$ws.Range("H132").Value = 2794.08
$ws.Range("I132").Value = 2613
$ws.Range("J132").Value = 3259.7144
$ws.Range("K132").Value = 7839
$ws.Range("L132").Value = 9779.143199999999
$ws.Range("M132").Value = -5309
$ws.Range("N132").Value = -14839.1432
